$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A,B,C -> B,C,D)
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold, border, centered) from the neighboring
# header cell onto the new header cell, then set its text.
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "ssim_dual"

# New ssim_dual data values for rows 2-51
$values = @(
    0.93520238,
    0.93691637000000005,
    0.85380898000000005,
    0.88221572000000004,
    0.85981969000000003,
    0.89952825999999997,
    0.92060065999999996,
    0.89254692999999996,
    0.90341241999999999,
    0.92406719000000004,
    0.92551298999999998,
    0.90270899000000004,
    0.87213996000000005,
    0.92796265,
    0.93648297000000003,
    0.89547111000000001,
    0.90644738000000002,
    0.90016898000000001,
    0.90088725999999997,
    0.90465205999999998,
    0.90041682999999995,
    0.85309853999999996,
    0.87858888000000002,
    0.90856221000000004,
    0.86770709999999995,
    0.88892360999999998,
    0.83904875999999995,
    0.91309035000000005,
    0.86627072000000005,
    0.90885209,
    0.83882851000000003,
    0.89761374000000005,
    0.85586609000000002,
    0.82659254999999998,
    0.89191284000000004,
    0.88055123999999996,
    0.90624934999999995,
    0.79560911999999995,
    0.89437677000000004,
    0.91013540000000004,
    0.83691934000000001,
    0.88430900999999995,
    0.86311130999999996,
    0.93718212999999995,
    0.98508467,
    0.87466685,
    0.92263298000000005,
    0.90042977000000002,
    0.92280094999999995,
    0.95688088999999998
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# New column A is narrower than the data columns (B:D keep their existing
# best-fit width of 12 automatically); match the target width of 11.
$ws.Columns.Item(1).ColumnWidth = 10.1

# Move the active selection to A2, as in the edited workbook
$ws.Range("A2").Select() | Out-Null
